$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column C header "13-01-2023", cloning B1 style (bold/border/centered).
# B1 itself ("06-01-2023") is left untouched by this edit.
$ws.Range("B1").Copy($ws.Range("C1")) | Out-Null
$ws.Range("C1").Value = '13-01-2023'

# Rewrite rows 2-45: fund rows now come first (alphabetical order), with the
# "avg" and "total" summary rows moved down to rows 44-45. Column B keeps the
# prior snapshot values realigned to the new row order, column C adds the new
# 13-01-2023 snapshot values.
$ws.Range("A2").Value = '1810 Renta variable'
$ws.Range("B2").Value = 65113.68
$ws.Range("C2").Value = 65208.4
$ws.Range("A3").Value = '1822 Raices Valores Negociables'
$ws.Range("B3").Value = 460076.59
$ws.Range("C3").Value = 432276.88
$ws.Range("A4").Value = 'Adcap IOL Acciones Argentina'
$ws.Range("B4").Value = 38928.01
$ws.Range("C4").Value = 49545.18
$ws.Range("A5").Value = 'Allaria Acciones'
$ws.Range("B5").Value = 36790.43
$ws.Range("C5").Value = 36852.77
$ws.Range("A6").Value = 'Alpha Acciones'
$ws.Range("B6").Value = 99318.73
$ws.Range("C6").Value = 99171.29
$ws.Range("A7").Value = 'Alpha Latam'
$ws.Range("B7").Value = 119.06
$ws.Range("C7").Value = 114.95
$ws.Range("A8").Value = 'Alpha Mega'
$ws.Range("B8").Value = 242684.39
$ws.Range("C8").Value = 243082.7
$ws.Range("A9").Value = 'Alpha Recursos Naturales'
$ws.Range("B9").Value = 220159.87
$ws.Range("C9").Value = 228455.66
$ws.Range("A10").Value = 'Alpha renta balan global'
$ws.Range("B10").Value = 179630.92
$ws.Range("C10").Value = 178136.02
$ws.Range("A11").Value = 'Argenfunds'
$ws.Range("B11").Value = 16550.17
$ws.Range("C11").Value = 16528.16
$ws.Range("A12").Value = 'Arpenta acciones'
$ws.Range("B12").Value = 1220.74
$ws.Range("C12").Value = 1214.03
$ws.Range("A13").Value = 'Arpenta ex Mercosur'
$ws.Range("B13").Value = 15021.03
$ws.Range("C13").Value = 15001.71
$ws.Range("A14").Value = 'Balanz'
$ws.Range("B14").Value = 342394.97
$ws.Range("C14").Value = 343069.22
$ws.Range("A15").Value = 'Bull Market'
$ws.Range("B15").Value = 145795.8
$ws.Range("C15").Value = 162896.96
$ws.Range("A16").Value = 'Compass Small Cap II'
$ws.Range("B16").Value = 80.71
$ws.Range("C16").Value = 82.77
$ws.Range("A17").Value = 'Consultatio Acciones Argentina'
$ws.Range("B17").Value = 350576.1
$ws.Range("C17").Value = 350043.47
$ws.Range("A18").Value = 'Consultatio Renta Variable'
$ws.Range("B18").Value = 125184.6
$ws.Range("C18").Value = 124903.09
$ws.Range("A19").Value = 'Delta Latinoamerica'
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("A20").Value = 'FBA Acciones Argentinas'
$ws.Range("B20").Value = 255598.38
$ws.Range("C20").Value = 255881.76
$ws.Range("A21").Value = 'FBA Calificado'
$ws.Range("B21").Value = 251419.77
$ws.Range("C21").Value = 251924.4
$ws.Range("A22").Value = 'Fima Acciones'
$ws.Range("B22").Value = 246329.97
$ws.Range("C22").Value = 246155.64
$ws.Range("A23").Value = 'Fima PB Acciones'
$ws.Range("B23").Value = 352513
$ws.Range("C23").Value = 369060.26
$ws.Range("A24").Value = 'Gainvest Renta Variable'
$ws.Range("B24").Value = 120873.18
$ws.Range("C24").Value = 120941.11
$ws.Range("A25").Value = 'Goal Acciones Argentinas'
$ws.Range("B25").Value = 50748.5
$ws.Range("C25").Value = 50811.18
$ws.Range("A26").Value = 'Goal acciones plus'
$ws.Range("B26").Value = 10488.52
$ws.Range("C26").Value = 10478.69
$ws.Range("A27").Value = 'HF Acciones Argentinas'
$ws.Range("B27").Value = 128297.94
$ws.Range("C27").Value = 128323.86
$ws.Range("A28").Value = 'HF Acciones Lideres'
$ws.Range("B28").Value = 174459.32
$ws.Range("C28").Value = 174240
$ws.Range("A29").Value = 'IAM Renta Variable'
$ws.Range("B29").Value = 53451.33
$ws.Range("C29").Value = 56055.15
$ws.Range("A30").Value = 'IEB Value'
$ws.Range("B30").Value = 2811.8
$ws.Range("C30").Value = 2833.91
$ws.Range("A31").Value = 'Lombardi'
$ws.Range("B31").Value = 62561.72
$ws.Range("C31").Value = 71093.2
$ws.Range("A32").Value = 'MAF'
$ws.Range("B32").Value = 13847.78
$ws.Range("C32").Value = 13990.63
$ws.Range("A33").Value = 'Megainver'
$ws.Range("B33").Value = 46332.05
$ws.Range("C33").Value = 46354.95
$ws.Range("A34").Value = 'Pellegrini Acciones'
$ws.Range("B34").Value = 49829.03
$ws.Range("C34").Value = 119916.14
$ws.Range("A35").Value = 'Pionero Acciones'
$ws.Range("B35").Value = 214583.84
$ws.Range("C35").Value = 214408.32
$ws.Range("A36").Value = 'Premier Renta Variable'
$ws.Range("B36").Value = 60645.8
$ws.Range("C36").Value = 60516.93
$ws.Range("A37").Value = 'Quinquela Acciones'
$ws.Range("B37").Value = 153737.81
$ws.Range("C37").Value = 154135.23
$ws.Range("A38").Value = 'Rofex 20 Renta Variable'
$ws.Range("B38").Value = 108208.8
$ws.Range("C38").Value = 108402.05
$ws.Range("A39").Value = 'Supefondo RV'
$ws.Range("B39").Value = 1986685.86
$ws.Range("C39").Value = 1988073.21
$ws.Range("A40").Value = 'Superfondo '
$ws.Range("B40").Value = 659915.41
$ws.Range("C40").Value = 661651.32
$ws.Range("A41").Value = 'Supergestion'
$ws.Range("B41").Value = 197661.1
$ws.Range("C41").Value = 197729.46
$ws.Range("A42").Value = 'Toronto Trust Multimercado'
$ws.Range("B42").Value = 79944.42
$ws.Range("C42").Value = 80133.28
$ws.Range("A43").Value = 'Toronto trust Argy'
$ws.Range("B43").Value = 123590.2
$ws.Range("C43").Value = 123591.93
$ws.Range("A44").Value = 'avg'
$ws.Range("B44").Value = 184385.27
$ws.Range("C44").Value = 186983
$ws.Range("A45").Value = 'total'
$ws.Range("B45").Value = 7744181.33
$ws.Range("C45").Value = 7853285.87
